$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requirements")

$ws.Range("B1").Value = "RELATIONS"
$ws.Range("B3").Value = "Parent: REQ-001"
$ws.Range("B4").Value = "Parent: REQ-002"
$ws.Range("B5").Value = "Parent: REQ-002`n----`nFile: src/samplefile.cpp"

$ws.Range("B1:B5").WrapText = $false
$ws.Columns.Item(2).ColumnWidth = 45.71

$wb.Worksheets.Item("Requirements").ListObjects.Item(1).ListColumns.Item(2).Name = "RELATIONS"
